# Add a new "signin_title" column (N) to the customer data sheet and
# populate it with "Sign In with Email" for every existing user row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "signin_title"

for ($row = 2; $row -le 4; $row++) {
    $ws.Cells.Item($row, 14).Value = "Sign In with Email"
}

$ws.Range("N4").Select()
